$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report previously listed a single short item (KENACOMB) in row 7, with
# the "total" row at 8 and the footer (timestamp / page / credit) at row 9.
# The new export lists four items (rows 7-10), pushing the total row down to
# 11 and the footer down to row 12.
# ---------------------------------------------------------------------------

# Insert three blank rows right before the old "total" row (row 8), shifting
# it (and the footer below it) down by three rows.
$ws.Rows.Item(8).Resize(3, 1).Insert()

# Clone row 7's cell formatting (styles / number formats / borders / fill)
# into the three freshly inserted rows so the new item rows look identical
# to the existing one.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Recreate the merged-cell layout for each new row (No / Name / Balance /
# Order-limit / Price columns), matching row 7's pattern.
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# Row heights matching the new export.
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5

# ---------------------------------------------------------------------------
# Helper: write a value as literal text (leading apostrophe) so numeric-
# looking strings such as "1", "36.00" or time-like "6:3" are stored as
# shared-string text instead of being coerced into numbers/dates.
# ---------------------------------------------------------------------------

# Row 7 - DECLOPHEN 75MG/3ML 3 AMPOULES
$ws.Range("A7").Value2 = 1
$ws.Range("C7").Value2 = "DECLOPHEN 75MG/3ML 3 AMPOULES"
$ws.Range("H7").Value2 = "'6:3"
$ws.Range("L7").Value2 = "'1"
$ws.Range("N7").Value2 = "'36.00"
$ws.Range("P7").Value2 = "'11.8800"
$ws.Range("Q7").Value2 = "'0:1"

# Row 8 - DICLAC 75 ID 30 TAB
$ws.Range("A8").Value2 = 2
$ws.Range("C8").Value2 = "DICLAC 75 ID 30 TAB"
$ws.Range("H8").Value2 = "'0:2"
$ws.Range("L8").Value2 = "'0"
$ws.Range("N8").Value2 = "'135.00"
$ws.Range("P8").Value2 = "'44.5500"
$ws.Range("Q8").Value2 = "'0:1"

# Row 9 - KENACOMB TOPICAL CREAM 15 GM
$ws.Range("A9").Value2 = 3
$ws.Range("C9").Value2 = "KENACOMB TOPICAL CREAM 15 GM"
$ws.Range("H9").Value2 = "'3:0"
$ws.Range("L9").Value2 = "'1"
$ws.Range("N9").Value2 = "'36.00"
$ws.Range("P9").Value2 = "'36.0000"
$ws.Range("Q9").Value2 = "'1:0"

# Row 10 - سرنجات 3 سم (3cc syringes)
$ws.Range("A10").Value2 = 4
$ws.Range("C10").Value2 = "سرنجات 3 سم"
$ws.Range("H10").Value2 = "'0:0"
$ws.Range("L10").Value2 = "'0"
$ws.Range("N10").Value2 = "'2.00"
$ws.Range("P10").Value2 = "'2.0000"
$ws.Range("Q10").Value2 = "'1:0"

# Total row (now row 11) - sum of the sale-price column for the 4 items.
$ws.Range("P11").Value2 = 94.43

# Footer row (now row 12) - updated generation timestamp.
$ws.Range("A12").Value2 = "Monday, 14 July, 2025 10:06 AM"

Write-Output "edit complete"
